# Capacity + Inv limit fix
# Changed code such that capacity values and investment values shown in
# interface are the actual values used instead of placeholder values.
#
# The Methanol_Plant row (row 7) previously only carried a placeholder
# investment-cost value of 1 in column F (2050). The other year columns
# (B=2020, C=2025, D=2030, E=2040) were left blank. This sets them to the
# actual placeholder value of 1 that is used consistently by the model
# (see note in column I / sharedStrings: "The methanol plant value is set
# to 1 to give the model a cost...").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investment_Cost")

$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1

# Reset the view back to the top-left corner / A1 so the sheet no longer
# opens scrolled down with a stray selection sitting on C26.
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
